$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    4289.535569812224,
    4229.735867409958,
    4191.67548425361,
    4191.087071822973,
    4191.087071822973,
    4176.000855543626,
    4107.624568105235,
    4026.648494151807,
    4008.77341602358,
    4008.77341602358,
    4008.77341602358
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
